$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column B/C updates (coin name/link reordering) ---
$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("B38").Value = "Hedera"
$ws.Range("C38").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("B46").Value = "Decentraland"
$ws.Range("C46").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("B48").Value = "RenderToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("B49").Value = "NEARProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"

# --- Column D (Price) updates ---
$ws.Range("D2").Value = "29.639.95"
$ws.Range("D3").Value = "1.943.25"
$ws.Range("D4").Value = "'0.9943"
$ws.Range("D5").Value = "'341.96"
$ws.Range("D6").Value = "'0.9930"
$ws.Range("D7").Value = "'0.4776"
$ws.Range("D8").Value = "'0.4103"
$ws.Range("D9").Value = "'48.61"
$ws.Range("D10").Value = "'0.08221"
$ws.Range("D11").Value = "'1.048"
$ws.Range("D12").Value = "'22.71"
$ws.Range("D13").Value = "1.917.67"
$ws.Range("D14").Value = "'6.139"
$ws.Range("D15").Value = "'7.442"
$ws.Range("D16").Value = "'92.60"
$ws.Range("D17").Value = "'0.9909"
$ws.Range("D18").Value = "'0.00001066"
$ws.Range("D19").Value = "'0.06663"
$ws.Range("D20").Value = "'18.07"
$ws.Range("D21").Value = "'0.9978"
$ws.Range("D22").Value = "29.641.15"
$ws.Range("D23").Value = "'5.624"
$ws.Range("D24").Value = "'11.29"
$ws.Range("D25").Value = "'2.253"
$ws.Range("D26").Value = "2.149.24"
$ws.Range("D27").Value = "'161.18"
$ws.Range("D28").Value = "'20.22"
$ws.Range("D29").Value = "'2.221"
$ws.Range("D30").Value = "'5.638"
$ws.Range("D31").Value = "'122.29"
$ws.Range("D32").Value = "'1.025"
$ws.Range("D33").Value = "'0.09646"
$ws.Range("D34").Value = "'1.470"
$ws.Range("D35").Value = "'3.630"
$ws.Range("D36").Value = "'5.497"
$ws.Range("D37").Value = "'0.02310"
$ws.Range("D38").Value = "'0.06228"
$ws.Range("D39").Value = "'8.698"
$ws.Range("D40").Value = "'1.207"
$ws.Range("D41").Value = "'0.6112"
$ws.Range("D42").Value = "'10.68"
$ws.Range("D43").Value = "'0.1915"
$ws.Range("D44").Value = "'0.9921"
$ws.Range("D45").Value = "'1.263"
$ws.Range("D46").Value = "'0.5711"
$ws.Range("D47").Value = "'12.51"
$ws.Range("D48").Value = "'2.321"
$ws.Range("D49").Value = "'2.006"
$ws.Range("D50").Value = "'0.07253"
$ws.Range("D51").Value = "'113.71"

# --- Column E (Volume 1h) updates ---
$ws.Range("E2").Value = "  +7.98%  "
$ws.Range("E3").Value = "  +6.28%  "
$ws.Range("E4").Value = "  -0.72%  "
$ws.Range("E5").Value = "  +3.07%  "
$ws.Range("E6").Value = "  -0.88%  "
$ws.Range("E7").Value = "  +4.47%  "
$ws.Range("E8").Value = "  +7.52%  "
$ws.Range("E9").Value = "  +4.86%  "
$ws.Range("E10").Value = "  +4.38%  "
$ws.Range("E11").Value = "  +8.64%  "
$ws.Range("E12").Value = "  +8.16%  "
$ws.Range("E13").Value = "  +4.38%  "
$ws.Range("E14").Value = "  +4.56%  "
$ws.Range("E15").Value = "  +5.11%  "
$ws.Range("E16").Value = "  +3.53%  "
$ws.Range("E17").Value = "  -1.14%  "
$ws.Range("E18").Value = "  +4.08%  "
$ws.Range("E19").Value = "  +0.93%  "
$ws.Range("E20").Value = "  +5.87%  "
$ws.Range("E21").Value = "  -0.29%  "
$ws.Range("E22").Value = "  +8.02%  "
$ws.Range("E23").Value = "  +5.60%  "
$ws.Range("E24").Value = "  +4.55%  "
$ws.Range("E25").Value = "  -1.32%  "
$ws.Range("E26").Value = "  +4.71%  "
$ws.Range("E27").Value = "  +3.55%  "
$ws.Range("E28").Value = "  +4.59%  "
$ws.Range("E29").Value = "  +7.41%  "
$ws.Range("E30").Value = "  +6.49%  "
$ws.Range("E31").Value = "  +3.39%  "
$ws.Range("E32").Value = "  +9.36%  "
$ws.Range("E33").Value = "  +3.81%  "
$ws.Range("E34").Value = "  +10.97%  "
$ws.Range("E35").Value = "  +1.50%  "
$ws.Range("E36").Value = "  +5.04%  "
$ws.Range("E37").Value = "  +6.15%  "
$ws.Range("E38").Value = "  +5.20%  "
$ws.Range("E39").Value = "  +7.28%  "
$ws.Range("E40").Value = "  +5.39%  "
$ws.Range("E41").Value = "  +5.92%  "
$ws.Range("E42").Value = "  +7.08%  "
$ws.Range("E43").Value = "  +5.10%  "
$ws.Range("E44").Value = "  -0.91%  "
$ws.Range("E45").Value = "  -0.27%  "
$ws.Range("E46").Value = "  +5.05%  "
$ws.Range("E47").Value = "  +5.09%  "
$ws.Range("E48").Value = "  +28.14%  "
$ws.Range("E49").Value = "  +7.44%  "
$ws.Range("E50").Value = "  +10.18%  "
$ws.Range("E51").Value = "  +2.98%  "
